# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 255-256) into the daily price log,
# pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 255 (existing rows 255.. shift to 257..)
$ws.Rows("255:256").Insert()

# --- New row 255: Zafiro rojo, Región de Arica y Parinacota ---
$ws.Cells.Item(255, 1).Value = 7
$ws.Cells.Item(255, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(255, 3).Value = "Ñuble"
$ws.Cells.Item(255, 4).Value2 = 44875
$ws.Cells.Item(255, 5).Value = 16
$ws.Cells.Item(255, 6).Value = 100112002
$ws.Cells.Item(255, 7).Value = "Pimiento"
$ws.Cells.Item(255, 8).Value = "Zafiro rojo"
$ws.Cells.Item(255, 9).Value = "Primera"
$ws.Cells.Item(255, 10).Value = 120
$ws.Cells.Item(255, 11).Value = 17000
$ws.Cells.Item(255, 12).Value = 18000
$ws.Cells.Item(255, 13).Value = 17500
$ws.Cells.Item(255, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(255, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(255, 16).Value = 1167
$ws.Cells.Item(255, 17).Value = 15
$ws.Cells.Item(255, 18).Value = "Hortaliza"

# --- New row 256: Zafiro verde, Región de Arica y Parinacota ---
$ws.Cells.Item(256, 1).Value = 7
$ws.Cells.Item(256, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(256, 3).Value = "Ñuble"
$ws.Cells.Item(256, 4).Value2 = 44875
$ws.Cells.Item(256, 5).Value = 16
$ws.Cells.Item(256, 6).Value = 100112002
$ws.Cells.Item(256, 7).Value = "Pimiento"
$ws.Cells.Item(256, 8).Value = "Zafiro verde"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 120
$ws.Cells.Item(256, 11).Value = 15000
$ws.Cells.Item(256, 12).Value = 16000
$ws.Cells.Item(256, 13).Value = 15500
$ws.Cells.Item(256, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(256, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(256, 16).Value = 1033
$ws.Cells.Item(256, 17).Value = 15
$ws.Cells.Item(256, 18).Value = "Hortaliza"
